$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting existing rows 117-217 down to 118-218
$ws.Rows(117).EntireRow.Insert()

# Populate the newly inserted row 117 with its data
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 44658
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = 100112043
$ws.Range("G117").Value = "Pepino ensalada"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 100
$ws.Range("K117").Value = 15000
$ws.Range("L117").Value = 16000
$ws.Range("M117").Value = 15500
$ws.Range("N117").Value = "$/caja 80 unidades"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 194
$ws.Range("Q117").Value = 80
$ws.Range("R117").Value = "Hortaliza"
